# Update the "想去人数" (interested-attendee) counts to the freshly
# generated values from the 456a3b4 output run.
# Two sheets carry the same event rows ("展览" and "全部类型"); both need
# the same F-column bump.

$wb = $excel.ActiveWorkbook

$exhibitionSheet = $wb.Worksheets.Item("展览")
$exhibitionSheet.Range("F2").Value = 119
$exhibitionSheet.Range("F3").Value = 2137
$exhibitionSheet.Range("F5").Value = 11121
$exhibitionSheet.Range("F10").Value = 11015
$exhibitionSheet.Range("F13").Value = 36
$exhibitionSheet.Range("F14").Value = 1714
$exhibitionSheet.Range("F15").Value = 5523
$exhibitionSheet.Range("F16").Value = 89

$allTypesSheet = $wb.Worksheets.Item("全部类型")
$allTypesSheet.Range("F2").Value = 119
$allTypesSheet.Range("F3").Value = 2137
$allTypesSheet.Range("F7").Value = 11121
$allTypesSheet.Range("F12").Value = 11015
$allTypesSheet.Range("F15").Value = 36
$allTypesSheet.Range("F16").Value = 1714
$allTypesSheet.Range("F17").Value = 5523
$allTypesSheet.Range("F18").Value = 89
